# Update workbook via Shiny-style data refresh
$wb = $excel.ActiveWorkbook

# --- Sheet "pro": new production values (B2:B26) ---
$proVals = @(2813111.8068729583,3287367.1300997585,3662740.8201511893,3227953.87468038,3312469.0678267083,2956003.0162836676,2856163.397867981,2850467.822401154,2604093.326811797,2535447.530723418,2798172.650688127,2611642.9799922295,2538626.8425992075,2634261.872707831,2829142.036533057,3051374.972308469,3202960.7410026803,3610021.6490008826,4151259.4542214484,4342123,4198234,4327768.984950564,4576163.316987814,4846140.944846592,4993665.023576268)
$wsPro = $wb.Worksheets.Item("pro")
for ($i = 0; $i -lt $proVals.Length; $i++) {
    $wsPro.Cells.Item(2 + $i, 2).Value = $proVals[$i]
}

# --- Sheet "ind": new index values (B2:B101) ---
$indVals = @(673991.9400906917,718269.7328207466,537534.72274806,522163.9170341421,787618.5171144058,839360.9896526561,628156.326574709,610194.1960459363,877554.0969938268,935204.8730711708,699883.4409808036,679870.2735732977,594625.2026426053,633689.01024709,474236.670246372,460675.8723810612,820401.8296278283,874297.996639456,654302.2902791037,635592.5159027878,759411.8555701455,809301.2959886121,605660.4195898242,588341.5595345482,720395.9076596692,767722.2016189747,574543.6873333694,558114.6102553135,687663.8361698807,732839.801926663,548438.5903334473,532755.988519496,628226.9853343992,669498.2566797521,501035.3956742215,486708.2882395003,611666.4646001904,651849.7951840076,487827.7378493687,473878.3033979367,756038.2212431404,805706.031293953,602969.8153188355,585727.8930942748,748005.6552520082,797145.767164475,596563.5322816558,579504.7977243206,662264.6951111678,705772.0683575297,528181.7898660754,513078.43127343233,750358.3300452789,799653.000552087,598439.8817118092,581327.4929414125,760520.8739147612,810483.1711160404,606544.904748972,589200.7528400487,845561.2418789255,901110.2524552912,674367.897721053,655084.3472881039,887474.4576990764,945776.9502866437,707795.313547914,687555.9061396864,955342.8560520076,1018103.9522169227,761922.9944918269,740135.8060150937,983492.9815293726,1048103.395686197,784373.8117697936,761944.6421597599,1099515.2516971482,1171747.7303401865,876906.0738304409,851830.9441322247,996468.7917014671,1061931.64963266,794722.5238365134,771997.43279883,1173486.0635956496,1250578.042898109,935900.6663473179,909138.5862412978,1249206.7156522898,1331273.1510835534,996290.8242831702,967801.8875753676,1255259.754272288,1337723.8431077853,1001118.3574372702,972491.3774160795,1334414.932188476,1422079.1077799248,1064247.683003873,1033815.5198808836)
$wsInd = $wb.Worksheets.Item("ind")
for ($i = 0; $i -lt $indVals.Length; $i++) {
    $wsInd.Cells.Item(2 + $i, 2).Value = $indVals[$i]
}

# --- Sheet "conso": new consumption values (B2:B26) ---
$consoVals = @(699386.9533236737,815838.8721687681,910544.6074912649,802216.676293145,821233.5910007586,733988.4185683706,708529.9949372734,705777.2005003267,642461.8318084963,622564.8424133114,686706.0278205522,641265.8194792488,623631.2257286278,646897.487138226,695326.3494886777,747563.1102776336,785418.8674807512,886298.068792961,1022894.6740994203,1067530,1142948,1274179.5840800432,1262130.9855690033,1336592.2986663848,1377280.250119308)
$wsConso = $wb.Worksheets.Item("conso")
for ($i = 0; $i -lt $consoVals.Length; $i++) {
    $wsConso.Cells.Item(2 + $i, 2).Value = $consoVals[$i]
}

# "VA" sheet (B2:B26) holds formulas `=pro!Bn - conso!Bn`; it recalculates
# automatically once the two input sheets above are updated.

# --- Column B on "pro" narrows slightly after the bestFit re-measure ---
$wsPro.Columns.Item(2).ColumnWidth = 13.76

# --- Refresh selections on every sheet to B2:B26, scroll back to top on "ind" ---
$wsPro.Range("B2:B26").Select()
$wsInd.Range("B2:B26").Select()
$wb.Worksheets.Item("VA").Range("B2:B26").Select()
$wsConso.Range("B2:B26").Select()

# Leave "pro" (the originally active sheet) selected/active at the end
$wsPro.Activate()
$wsPro.Range("B2:B26").Select()
